$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9
$ws.Range("G2").Value = 5
$ws.Range("J2").Value = 8

# Row 3
$ws.Range("D3").Value = 14
$ws.Range("G3").Value = 8
$ws.Range("J3").Value = 12
$ws.Range("M3").Value = 9

# Row 4
$ws.Range("D4").Value = 18
$ws.Range("G4").Value = 10
$ws.Range("J4").Value = 15
$ws.Range("M4").Value = 11

# Row 5
$ws.Range("D5").Value = 21
$ws.Range("G5").Value = 12
$ws.Range("J5").Value = 17
$ws.Range("M5").Value = 13

# Row 6
$ws.Range("D6").Value = 22
$ws.Range("G6").Value = 13
$ws.Range("J6").Value = 18
$ws.Range("M6").Value = 15

# Row 7
$ws.Range("D7").Value = 23
$ws.Range("G7").Value = 14
$ws.Range("J7").Value = 19
$ws.Range("M7").Value = 17

# Row 8
$ws.Range("M8").Value = 13

# Row 9
$ws.Range("M9").Value = 11

# Row 10
$ws.Range("P10").Value = 14

# Row 11
$ws.Range("P11").Value = 12

# Update selection to reflect the saved view state (row 5 selected)
$ws.Rows(5).Select()
